$d = $word.ActiveDocument
$tab = [char]9
$cr = [char]13

# ---------------------------------------------------------------------
# 1) Merge the two runs "SAT Feb 24" + " 10:24:26 PST 2018" into a
#    single run "SAT Feb 24 10:24:26 PST 2018".
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("SAT Feb 24 10:24:26 PST 2018", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "SAT Feb 24 10:24:26 PST 2018", 2)

# ---------------------------------------------------------------------
# 2) Append a new purchase-details block (MAMATHA / RAMANNA, 06/03/2018)
#    right after the last "Amount Received mode ... - CASH" paragraph,
#    which sits just before the long run of trailing blank paragraphs.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$idx = -1
for ($i = $count; $i -ge 1; $i--) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -eq ("Amount Received mode" + $tab + $tab + "- CASH" + $cr)) {
        $idx = $i
        break
    }
}

$script:cur = $idx

function InsertNewPara([string]$text) {
    $r = $d.Paragraphs($script:cur).Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $script:cur = $script:cur + 1
    if ($text -ne $null -and $text -ne "") {
        $newRange = $d.Paragraphs($script:cur).Range
        $newRange.InsertBefore($text)
    }
}

# blank separator line
InsertNewPara ""

# timestamp line
InsertNewPara ("MON Mar 05" + " 10:50:15 IST 2018")

# Person Name
InsertNewPara ("Person Name" + $tab + $tab + $tab + $tab + "- RAMANNA")

# Bill number
InsertNewPara ("Bill number" + $tab + $tab + $tab + $tab + "- 3701")

# dashed separator
InsertNewPara ("---------------------------------------------------------------")

# Item Name
InsertNewPara ("Item Name" + $tab + $tab + $tab + $tab + "- POTATO")

# Number of Pockets
InsertNewPara ("Number of Pockets" + $tab + $tab + $tab + "- 1")

# Number of KGs
InsertNewPara ("Number of KGs" + $tab + $tab + $tab + "- 48")

# Rate
InsertNewPara ("Rate" + $tab + $tab + $tab + $tab + $tab + "- 14")

# Total Price
InsertNewPara ("Total Price" + $tab + $tab + $tab + $tab + "- 672.0")

# Amount balance (bold)
$r = $d.Paragraphs($script:cur).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$script:cur = $script:cur + 1
$newRange = $d.Paragraphs($script:cur).Range
$newRange.InsertBefore("Amount balance" + $tab + $tab + $tab + "- 1348.0")
$newRange.Bold = 1

# trailing blank paragraph (new)
InsertNewPara ""

Write-Output "done. paragraphs=$($d.Paragraphs.Count)"
